$d = $word.ActiveDocument

# Remove the empty footnote (id 22, containing just a single Tibetan
# punctuation mark) together with its in-text reference mark. Deleting
# via the Footnotes collection removes both the footnoteReference run
# in the body and the w:footnote definition in footnotes.xml.
$fn = $d.Footnotes.Item($d.Footnotes.Count)
$fn.Delete()
